$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.05463385581970215
$ws.Range("B2").Value = 0.9848271608352661
$ws.Range("C2").Value = 0.008333253674209118
$ws.Range("D2").Value = 0.9983949661254883
$ws.Range("A3").Value = 0.009050683118402958
$ws.Range("B3").Value = 0.9983867406845093
$ws.Range("C3").Value = 0.006796804256737232
$ws.Range("D3").Value = 0.9987159967422485
$ws.Range("A4").Value = 0.005266681313514709
$ws.Range("B4").Value = 0.9985501170158386
$ws.Range("C4").Value = 0.001667353790253401
$ws.Range("D4").Value = 0.9995987415313721
$ws.Range("A5").Value = 0.002695761388167739
$ws.Range("B5").Value = 0.9993261098861694
$ws.Range("C5").Value = 0.0003038088616449386
$ws.Range("D5").Value = 1
$ws.Range("A6").Value = 0.001647535827942193
$ws.Range("B6").Value = 0.9995711445808411
$ws.Range("C6").Value = 0.000277533516054973
$ws.Range("D6").Value = 0.9998394846916199
$ws.Range("A7").Value = 0.001343438751064241
$ws.Range("B7").Value = 0.9996936917304993
$ws.Range("C7").Value = 0.0005422402173280716
$ws.Range("D7").Value = 0.9996789693832397
$ws.Range("A8").Value = 0.00114441174082458
$ws.Range("B8").Value = 0.9996732473373413
$ws.Range("C8").Value = 0.0004783542244695127
$ws.Range("D8").Value = 0.9998394846916199
$ws.Range("A9").Value = 0.001100060646422207
$ws.Range("B9").Value = 0.9996732473373413
$ws.Range("C9").Value = 0.0004847114905714989
$ws.Range("D9").Value = 0.9996789693832397
$ws.Range("A10").Value = 0.0009437543340027332
$ws.Range("B10").Value = 0.9998570680618286
$ws.Range("C10").Value = [double]"7.103585812728852E-05"
$ws.Range("D10").Value = 1
$ws.Range("A11").Value = 0.0007627906743437052
$ws.Range("B11").Value = 0.9997753500938416
$ws.Range("C11").Value = 0.001276141148991883
$ws.Range("D11").Value = 0.9997592568397522
$ws.Range("A12").Value = 0.001428489573299885
$ws.Range("B12").Value = 0.9997345209121704
$ws.Range("C12").Value = [double]"2.141471850336529E-05"
$ws.Range("D12").Value = 1
$ws.Range("A13").Value = 0.0003187571419402957
$ws.Range("B13").Value = 0.9998978972434998
$ws.Range("C13").Value = [double]"2.418840267637279E-05"
$ws.Range("D13").Value = 1
$ws.Range("A14").Value = 0.0001668044278630987
$ws.Range("B14").Value = 0.9999387264251709
$ws.Range("C14").Value = [double]"7.552432361990213E-05"
$ws.Range("D14").Value = 0.9999197721481323
$ws.Range("A15").Value = 0.0002086235617753118
$ws.Range("B15").Value = 0.9999591708183289
$ws.Range("C15").Value = [double]"6.022910019964911E-05"
$ws.Range("D15").Value = 1
$ws.Range("A16").Value = 0.0008502230048179626
$ws.Range("B16").Value = 0.9997957944869995
$ws.Range("C16").Value = [double]"6.3439438235946E-06"
$ws.Range("D16").Value = 1
$ws.Range("A17").Value = 0.0009410029160790145
$ws.Range("B17").Value = 0.9997957944869995
$ws.Range("C17").Value = [double]"5.832837814523373E-06"
$ws.Range("D17").Value = 1
$ws.Range("A18").Value = 0.000374061957700178
$ws.Range("B18").Value = 0.9998774528503418
$ws.Range("C18").Value = [double]"3.66475524060661E-06"
$ws.Range("D18").Value = 1
$ws.Range("A19").Value = 0.0005728129763156176
$ws.Range("B19").Value = 0.9998774528503418
$ws.Range("C19").Value = [double]"7.164574981288752E-06"
$ws.Range("D19").Value = 1
$ws.Range("A20").Value = 0.0002683688944671303
$ws.Range("B20").Value = 0.9998978972434998
$ws.Range("C20").Value = [double]"2.085754567815457E-05"
$ws.Range("D20").Value = 1
$ws.Range("A21").Value = 0.0004725077305920422
$ws.Range("B21").Value = 0.9998774528503418
$ws.Range("C21").Value = [double]"1.000747033685911E-05"
$ws.Range("D21").Value = 1
$ws.Range("A22").Value = 0.0002783602685667574
$ws.Range("B22").Value = 0.9998570680618286
$ws.Range("C22").Value = [double]"8.233239350374788E-05"
$ws.Range("D22").Value = 0.9999197721481323
$ws.Range("A23").Value = [double]"3.40571059496142E-05"
$ws.Range("B23").Value = 0.999979555606842
$ws.Range("C23").Value = [double]"1.646588202675048E-06"
$ws.Range("D23").Value = 1
$ws.Range("A24").Value = [double]"4.139623342780396E-06"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = [double]"1.726305868032796E-06"
$ws.Range("D24").Value = 1
$ws.Range("A25").Value = 0.0001036266112350859
$ws.Range("B25").Value = 0.9999591708183289
$ws.Range("C25").Value = [double]"4.269593773642555E-06"
$ws.Range("D25").Value = 1
$ws.Range("A26").Value = 0.0004910431453026831
$ws.Range("B26").Value = 0.9997957944869995
$ws.Range("C26").Value = [double]"9.000656973512378E-06"
$ws.Range("D26").Value = 1
$ws.Range("A27").Value = 0.0003598150797188282
$ws.Range("B27").Value = 0.9999183416366577
$ws.Range("C27").Value = [double]"4.982568952982547E-07"
$ws.Range("D27").Value = 1
$ws.Range("A28").Value = [double]"4.669967893278226E-05"
$ws.Range("B28").Value = 0.999979555606842
$ws.Range("C28").Value = [double]"4.337390180353395E-07"
$ws.Range("D28").Value = 1
$ws.Range("A29").Value = [double]"1.625025288376492E-05"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = [double]"7.744901085970923E-05"
$ws.Range("D29").Value = 1
$ws.Range("A30").Value = 0.000184772870852612
$ws.Range("B30").Value = 0.9999591708183289
$ws.Range("C30").Value = [double]"4.94184450872126E-06"
$ws.Range("D30").Value = 1
$ws.Range("A31").Value = 0.001243882812559605
$ws.Range("B31").Value = 0.9998162388801575
$ws.Range("C31").Value = [double]"3.348319296492264E-05"
$ws.Range("D31").Value = 1
$ws.Range("A32").Value = 0.0003245847474317998
$ws.Range("B32").Value = 0.9999183416366577
$ws.Range("C32").Value = [double]"7.342269441323879E-07"
$ws.Range("D32").Value = 1
$ws.Range("A33").Value = [double]"5.498053360497579E-05"
$ws.Range("B33").Value = 0.999979555606842
$ws.Range("C33").Value = [double]"3.527142382608872E-07"
$ws.Range("D33").Value = 1
$ws.Range("A34").Value = 0.000274541846010834
$ws.Range("B34").Value = 0.9999183416366577
$ws.Range("C34").Value = [double]"2.451183718221728E-05"
$ws.Range("D34").Value = 1
$ws.Range("A35").Value = 0.000220989910303615
$ws.Range("B35").Value = 0.999979555606842
$ws.Range("C35").Value = [double]"6.711389687552582E-06"
$ws.Range("D35").Value = 1
$ws.Range("A36").Value = 0.0005764481029473245
$ws.Range("B36").Value = 0.9999183416366577
$ws.Range("C36").Value = [double]"2.90350362774916E-05"
$ws.Range("D36").Value = 1
$ws.Range("A37").Value = 0.0001323032629443333
$ws.Range("B37").Value = 0.9999387264251709
$ws.Range("C37").Value = [double]"1.301797510677716E-06"
$ws.Range("D37").Value = 1
$ws.Range("A38").Value = 0.0004009020631201565
$ws.Range("B38").Value = 0.9998978972434998
$ws.Range("C38").Value = [double]"1.390770648868056E-05"
$ws.Range("D38").Value = 1
$ws.Range("A39").Value = [double]"5.005372440791689E-05"
$ws.Range("B39").Value = 0.999979555606842
$ws.Range("C39").Value = [double]"2.648499958013417E-06"
$ws.Range("D39").Value = 1
$ws.Range("A40").Value = 0.0002380369114689529
$ws.Range("B40").Value = 0.9999183416366577
$ws.Range("C40").Value = [double]"1.938965255021685E-07"
$ws.Range("D40").Value = 1
$ws.Range("A41").Value = 0.0006746638100594282
$ws.Range("B41").Value = 0.9998570680618286
$ws.Range("C41").Value = [double]"5.487586895469576E-05"
$ws.Range("D41").Value = 1
$ws.Range("A42").Value = 0.0001714193786028773
$ws.Range("B42").Value = 0.9999591708183289
$ws.Range("C42").Value = [double]"1.916809196700342E-06"
$ws.Range("D42").Value = 1
$ws.Range("A43").Value = [double]"5.287936346576316E-06"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = [double]"7.07420326762076E-07"
$ws.Range("D43").Value = 1
$ws.Range("A44").Value = [double]"3.981957434007199E-06"
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = [double]"1.474851387683884E-07"
$ws.Range("D44").Value = 1
$ws.Range("A45").Value = 0.0009792785858735442
$ws.Range("B45").Value = 0.9998774528503418
$ws.Range("C45").Value = [double]"1.053171786224993E-06"
$ws.Range("D45").Value = 1
$ws.Range("A46").Value = 0.0002712005807552487
$ws.Range("B46").Value = 0.9999387264251709
$ws.Range("C46").Value = [double]"4.352221196768369E-07"
$ws.Range("D46").Value = 1
$ws.Range("A47").Value = [double]"3.098827437497675E-05"
$ws.Range("B47").Value = 0.999979555606842
$ws.Range("C47").Value = [double]"3.474413290405209E-07"
$ws.Range("D47").Value = 1
$ws.Range("A48").Value = [double]"1.451865045964951E-05"
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = [double]"6.213860359594037E-08"
$ws.Range("D48").Value = 1
$ws.Range("A49").Value = 0.0002881486434489489
$ws.Range("B49").Value = 0.9999387264251709
$ws.Range("C49").Value = [double]"2.998918171215337E-05"
$ws.Range("D49").Value = 1
$ws.Range("A50").Value = 0.0001554638729430735
$ws.Range("B50").Value = 0.9999387264251709
$ws.Range("C50").Value = 0.0001027780162985437
$ws.Range("D50").Value = 1
$ws.Range("A51").Value = 0.0005663592019118369
$ws.Range("B51").Value = 0.9998774528503418
$ws.Range("C51").Value = [double]"8.347932123342616E-08"
$ws.Range("D51").Value = 1
